$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2204.2292171047425
$ws.Range("B1").Value = 1384.3421185015861
$ws.Range("C1").Value = 1431.1789760132535
$ws.Range("A2").Value = 2227.1192368922416
$ws.Range("B2").Value = 1482.0442411429403
$ws.Range("C2").Value = 1329.6419960545563
$ws.Range("A3").Value = 2343.913706953103
$ws.Range("B3").Value = 1585.8742932348805
$ws.Range("C3").Value = 1454.0924806953633
$ws.Range("A4").Value = 2320.9625221260826
$ws.Range("B4").Value = 1774.3335864583587
$ws.Range("C4").Value = 1704.2599074406442
$ws.Range("A5").Value = 2422.6539934395523
$ws.Range("B5").Value = 1668.399809907037
$ws.Range("C5").Value = 1626.0505420097888
$ws.Range("A6").Value = 2360.708401719302
$ws.Range("B6").Value = 1774.9878092462952
$ws.Range("C6").Value = 1784.9998764703353
$ws.Range("A7").Value = 1992.9848472506753
$ws.Range("B7").Value = 1566.2947138691036
$ws.Range("C7").Value = 1482.2947108097806
$ws.Range("A8").Value = 2135.7478628374106
$ws.Range("B8").Value = 1655.7390014213242
$ws.Range("C8").Value = 1507.3585028172008
$ws.Range("A9").Value = 2471.0153725333803
$ws.Range("B9").Value = 1788.3101470495083
$ws.Range("C9").Value = 1513.896066139911
$ws.Range("A10").Value = 2111.594404266449
$ws.Range("B10").Value = 1365.1595641833082
$ws.Range("C10").Value = 1300.9352233413354
$ws.Range("A11").Value = 1970.003883917462
$ws.Range("B11").Value = 1416.2529907667447
$ws.Range("C11").Value = 1298.4013275002937
$ws.Range("A12").Value = 2787.95371453083
$ws.Range("B12").Value = 2270.4399898165884
$ws.Range("C12").Value = 2036.7367177336491
$ws.Range("A13").Value = 2315.185341962035
$ws.Range("B13").Value = 1782.3539396052565
$ws.Range("C13").Value = 1791.9062670150788
$ws.Range("A14").Value = 2593.0481324833722
$ws.Range("B14").Value = 1922.4491827528125
$ws.Range("C14").Value = 1703.8048925083233
$ws.Range("A15").Value = 2508.510880379053
$ws.Range("B15").Value = 2026.7228915761436
$ws.Range("C15").Value = 1827.5083118642979
$ws.Range("A16").Value = 2205.521116644954
$ws.Range("B16").Value = 1510.1697609818998
$ws.Range("C16").Value = 1270.5287368582376
$ws.Range("A17").Value = 2225.544813078163
$ws.Range("B17").Value = 1682.2602014089707
$ws.Range("C17").Value = 1573.665397592627
$ws.Range("A18").Value = 2487.907268443025
$ws.Range("B18").Value = 2061.2592446629073
$ws.Range("C18").Value = 1916.9784859400013
$ws.Range("A19").Value = 1742.2116489556101
$ws.Range("B19").Value = 1927.5181691925923
$ws.Range("C19").Value = 1884.7058103212241
$ws.Range("A20").Value = 2351.2568365319244
$ws.Range("B20").Value = 1850.0736616385993
$ws.Range("C20").Value = 1653.346311004787
$ws.Range("A21").Value = 2583.1100797798213
$ws.Range("B21").Value = 1907.0314517922743
$ws.Range("C21").Value = 1813.0616860349023
$ws.Range("A22").Value = 2448.3691808101894
$ws.Range("B22").Value = 1890.2457075819934
$ws.Range("C22").Value = 1644.281166023012
